$d = $word.ActiveDocument

$replacements = @(
    @("68-43=25", "75+8=83"),
    @("80-38=42", "57+26=83"),
    @("0+94=94", "26+54=80"),
    @("49+12=61", "29-22=7"),
    @("36+46=82", "87-64=23"),
    @("91-5=86", "21+36=57"),
    @("26+6=32", "68-6=62"),
    @("2+7=9", "60+19=79"),
    @("40-14=26", "42+25=67"),
    @("9+60=69", "19+34=53"),
    @("18+8=26", "95-75=20"),
    @("68-32=36", "0+61=61"),
    @("15+55=70", "3+51=54"),
    @("12+38=50", "5+87=92"),
    @("83-7=76", "98-89=9"),
    @("72-43=29", "1+45=46"),
    @("96-58=38", "24+11=35"),
    @("13+79=92", "33+62=95"),
    @("87-26=61", "7-7=0"),
    @("20+35=55", "2+24=26"),
    @("53-19=34", "33-13=20"),
    @("2+9=11", "40+38=78"),
    @("8+64=72", "78-35=43"),
    @("5+56=61", "91-72=19"),
    @("1+50=51", "93-55=38"),
    @("7+54=61", "88+7=95"),
    @("90-14=76", "73-8=65"),
    @("17+17=34", "54+39=93"),
    @("19-5=14", "10+37=47"),
    @("93-79=14", "29-13=16"),
    @("95-23=72", "43-4=39"),
    @("17+57=74", "70-62=8"),
    @("47-34=13", "71+26=97"),
    @("2+4=6", "98-10=88"),
    @("66+20=86", "76-8=68"),
    @("38+56=94", "61-47=14"),
    @("47+47=94", "38+2=40"),
    @("45+41=86", "89-46=43"),
    @("22+43=65", "74-46=28"),
    @("38+11=49", "33-10=23"),
    @("56-18=38", "52-30=22"),
    @("56+41=97", "99-93=6"),
    @("17+42=59", "64-3=61"),
    @("20+46=66", "69-66=3"),
    @("52-31=21", "23+9=32"),
    @("28+26=54", "97-17=80"),
    @("71-53=18", "47+9=56"),
    @("52-41=11", "69+1=70"),
    @("5+81=86", "3+68=71"),
    @("53-0=53", "30+41=71"),
    @("14+25=39", "74-19=55"),
    @("88-17=71", "50+36=86"),
    @("39+19=58", "26+7=33"),
    @("84-72=12", "23+76=99"),
    @("15+81=96", "85-18=67"),
    @("14+59=73", "93+6=99"),
    @("94-25=69", "91-87=4"),
    @("52-52=0", "18+66=84"),
    @("37+2=39", "72-66=6"),
    @("79-77=2", "62-51=11"),
    @("45+5=50", "62-13=49"),
    @("14+38=52", "37+9=46"),
    @("18-2=16", "83+13=96"),
    @("73-33=40", "16+49=65"),
    @("38-20=18", "30-14=16"),
    @("24+40=64", "19+45=64"),
    @("67-2=65", "82-55=27"),
    @("68+8=76", "5+3=8"),
    @("50+4=54", "73+22=95"),
    @("27+26=53", "90-11=79"),
    @("63-7=56", "50-16=34"),
    @("48+13=61", "0+9=9"),
    @("10+33=43", "25+69=94"),
    @("78-36=42", "3+88=91"),
    @("89-23=66", "30+28=58"),
    @("5+71=76", "41+18=59"),
    @("39+31=70", "25+16=41"),
    @("91-11=80", "28+19=47"),
    @("42-37=5", "1+60=61"),
    @("23+34=57", "51+24=75"),
    @("14+8=22", "21+51=72"),
    @("94-37=57", "78-65=13"),
    @("99-45=54", "84-37=47"),
    @("94-61=33", "14+5=19"),
    @("72-57=15", "34+31=65"),
    @("51-47=4", "45+20=65"),
    @("29+39=68", "8+25=33"),
    @("96-36=60", "87+0=87"),
    @("3+4=7", "85-10=75"),
    @("41+21=62", "2+86=88"),
    @("81-49=32", "70-9=61"),
    @("12+71=83", "87-75=12"),
    @("74+5=79", "69-3=66"),
    @("50-22=28", "17+34=51"),
    @("92-37=55", "28+22=50"),
    @("74-32=42", "27+12=39"),
    @("63+17=80", "86-12=74"),
    @("12+52=64", "34+61=95"),
    @("16-0=16", "69-35=34"),
    @("93-18=75", "5+2=7"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Done: applied $($replacements.Count) replacements"
